# Updates crypto price/volume figures (and swaps two coin rows / replaces
# Filecoin with OKB) to match the latest scraped snapshot.
# NOTE: a handful of "Price" cells below are written with a leading
# apostrophe. Those values are plain decimal numbers whose text
# representation must be preserved exactly (e.g. trailing zero such as
# "259.11"/"0.150", or a fixed-point form like "0.0000212"). Excel's
# automatic type detection would otherwise silently convert them into
# floating point numbers and mangle the text, so we force them to stay
# text cells the same way a user typing them in the UI would.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '98.623.84'
$ws.Range('E2').Value2 = '  -0.30%  '
$ws.Range('D3').Value2 = '3.341.88'
$ws.Range('E3').Value2 = '  -1.51%  '
$ws.Range('E4').Value2 = '  -0.03%  '
$ws.Range('D5').Value2 = '''259.11'
$ws.Range('E5').Value2 = '  -0.93%  '
$ws.Range('D6').Value2 = '''646.76'
$ws.Range('E6').Value2 = '  +1.81%  '
$ws.Range('D7').Value2 = '1.54'
$ws.Range('E7').Value2 = '  +9.84%  '
$ws.Range('E8').Value2 = '  +15.45%  '
$ws.Range('D9').Value2 = '''1.09'
$ws.Range('E9').Value2 = '  +23.72%  '
$ws.Range('E10').Value2 = '  -0.01%  '
$ws.Range('D11').Value2 = '3.338.16'
$ws.Range('E11').Value2 = '  -1.56%  '
$ws.Range('B12').Value2 = 'TRON'
$ws.Range('C12').Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value2 = '''0.208'
$ws.Range('E12').Value2 = '  +3.80%  '
$ws.Range('B13').Value2 = 'Avalanche'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value2 = '''43.94'
$ws.Range('E13').Value2 = '  +20.40%  '
$ws.Range('E14').Value2 = '  +7.30%  '
$ws.Range('D15').Value2 = '98.691.79'
$ws.Range('E15').Value2 = '  -0.03%  '
$ws.Range('D16').Value2 = '3.961.83'
$ws.Range('E16').Value2 = '  -1.09%  '
$ws.Range('D17').Value2 = '''5.56'
$ws.Range('E17').Value2 = '  -0.05%  '
$ws.Range('D18').Value2 = '3.331.18'
$ws.Range('E18').Value2 = '  -1.60%  '
$ws.Range('D19').Value2 = '''7.43'
$ws.Range('E19').Value2 = '  +19.16%  '
$ws.Range('E20').Value2 = '  +9.47%  '
$ws.Range('D21').Value2 = '''536.71'
$ws.Range('E21').Value2 = '  +7.81%  '
$ws.Range('E22').Value2 = '  -1.48%  '
$ws.Range('D23').Value2 = '''10.20'
$ws.Range('E23').Value2 = '  +7.95%  '
$ws.Range('D24').Value2 = '''0.0000212'
$ws.Range('E24').Value2 = '  -0.83%  '
$ws.Range('D25').Value2 = '0.434'
$ws.Range('E25').Value2 = '  +52.37%  '
$ws.Range('D26').Value2 = '103.31'
$ws.Range('E26').Value2 = '  +15.52%  '
$ws.Range('D27').Value2 = '''6.23'
$ws.Range('E27').Value2 = '  +8.24%  '
$ws.Range('D28').Value2 = '12.69'
$ws.Range('E28').Value2 = '  +4.14%  '
$ws.Range('D29').Value2 = '3.518.89'
$ws.Range('E29').Value2 = '  -1.14%  '
$ws.Range('D30').Value2 = '''0.150'
$ws.Range('E30').Value2 = '  +13.38%  '
$ws.Range('E31').Value2 = '  +0.46%  '
$ws.Range('E32').Value2 = '  +13.17%  '
$ws.Range('E33').Value2 = '  -7.28%  '
$ws.Range('E34').Value2 = '  -0.08%  '
$ws.Range('D35').Value2 = '''29.17'
$ws.Range('E35').Value2 = '  +3.82%  '
$ws.Range('D36').Value2 = '''0.527'
$ws.Range('E36').Value2 = '  +11.07%  '
$ws.Range('E37').Value2 = '  +4.21%  '
$ws.Range('E38').Value2 = '  +2.73%  '
$ws.Range('E39').Value2 = '  +3.36%  '
$ws.Range('D40').Value2 = '''515.01'
$ws.Range('E40').Value2 = '  +1.84%  '
$ws.Range('D41').Value2 = '''24.70'
$ws.Range('E41').Value2 = '  -0.66%  '
$ws.Range('D42').Value2 = '''1.32'
$ws.Range('E42').Value2 = '  +2.30%  '
$ws.Range('D43').Value2 = '3.77'
$ws.Range('E43').Value2 = '  +0.94%  '
$ws.Range('E44').Value2 = '  -2.18%  '
$ws.Range('D45').Value2 = '''0.809'
$ws.Range('E45').Value2 = '  +1.82%  '
$ws.Range('D46').Value2 = '''0.0400'
$ws.Range('E46').Value2 = '  +21.89%  '
$ws.Range('D48').Value2 = '''2.03'
$ws.Range('E48').Value2 = '  +3.15%  '
$ws.Range('D49').Value2 = '''163.73'
$ws.Range('E49').Value2 = '  +1.89%  '
$ws.Range('D50').Value2 = '''7.73'
$ws.Range('E50').Value2 = '  +16.59%  '
$ws.Range('B51').Value2 = 'OKB'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value2 = '''49.74'
$ws.Range('E51').Value2 = '  +6.68%  '
